$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "country" column (D) values were originally stored with a leading
# non-breaking space (e.g. "\u00A0USA"). Re-write each one without the
# leading non-breaking space so the shared-string table no longer carries
# the stray whitespace variants (Excel will naturally dedupe/re-pack the
# shared-string table on save).
$countries = @{
    2  = "USA"
    3  = "KOR"
    4  = "USA"
    5  = "USA"
    6  = "CAN"
    7  = "KOR"
    8  = "DNK"
    9  = "USA"
    10 = "USA"
    11 = "USA"
    12 = "USA"
    13 = "KOR"
    14 = "USA"
    15 = "KOR"
    16 = "SWE"
    17 = "CHN"
    18 = "KOR"
    19 = "KOR"
    20 = "UKR"
    21 = "SWE"
}

foreach ($row in $countries.Keys) {
    $ws.Cells.Item($row, 4).Value = $countries[$row]
}

# Move the active selection to D22, matching the saved view state.
$ws.Range("D22").Select()
